# Scheduled runner update: refresh computed market-price / profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for specific leve rows across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1544623.5
$ws.Cells.Item(137, 9).Value = 2689429
$ws.Cells.Item(137, 10).Value = 1624.7391
$ws.Cells.Item(137, 11).Value = 8068287
$ws.Cells.Item(137, 12).Value = 4874.2173
$ws.Cells.Item(137, 13).Value = -8065737
$ws.Cells.Item(137, 14).Value = -9974.2173


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 12824215
$ws.Cells.Item(61, 9).Value = 22224394
$ws.Cells.Item(61, 10).Value = 5790.909
$ws.Cells.Item(61, 11).Value = 22224394
$ws.Cells.Item(61, 12).Value = 5790.909
$ws.Cells.Item(61, 13).Value = -22224182
$ws.Cells.Item(61, 14).Value = -6214.909

$ws.Cells.Item(74, 8).Value = 29415110
$ws.Cells.Item(74, 9).Value = 2481
$ws.Cells.Item(74, 10).Value = 38465150
$ws.Cells.Item(74, 11).Value = 2481
$ws.Cells.Item(74, 12).Value = 38465150
$ws.Cells.Item(74, 13).Value = -1607
$ws.Cells.Item(74, 14).Value = -38466898

$ws.Cells.Item(77, 8).Value = 29415110
$ws.Cells.Item(77, 9).Value = 2481
$ws.Cells.Item(77, 10).Value = 38465150
$ws.Cells.Item(77, 11).Value = 12405
$ws.Cells.Item(77, 12).Value = 192325750
$ws.Cells.Item(77, 13).Value = -8037
$ws.Cells.Item(77, 14).Value = -192334486

$ws.Cells.Item(97, 8).Value = 1131.9048
$ws.Cells.Item(97, 9).Value = 1051.1111
$ws.Cells.Item(97, 10).Value = 1192.5
$ws.Cells.Item(97, 11).Value = 1051.1111
$ws.Cells.Item(97, 12).Value = 1192.5
$ws.Cells.Item(97, 13).Value = -555.1111000000001
$ws.Cells.Item(97, 14).Value = -2184.5

$ws.Cells.Item(136, 8).Value = 12824215
$ws.Cells.Item(136, 9).Value = 22224394
$ws.Cells.Item(136, 10).Value = 5790.909
$ws.Cells.Item(136, 11).Value = 66673182
$ws.Cells.Item(136, 12).Value = 17372.727
$ws.Cells.Item(136, 13).Value = -66670632
$ws.Cells.Item(136, 14).Value = -22472.727


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2252.5
$ws.Cells.Item(99, 9).Value = 1755
$ws.Cells.Item(99, 10).Value = 2750
$ws.Cells.Item(99, 11).Value = 1755
$ws.Cells.Item(99, 12).Value = 2750
$ws.Cells.Item(99, 13).Value = -257
$ws.Cells.Item(99, 14).Value = -5746

$ws.Cells.Item(132, 8).Value = 76158.125
$ws.Cells.Item(132, 10).Value = 75902
$ws.Cells.Item(132, 12).Value = 75902
$ws.Cells.Item(132, 14).Value = -86022

$ws.Cells.Item(134, 8).Value = 1996.6349
$ws.Cells.Item(134, 9).Value = 1840.6072
$ws.Cells.Item(134, 11).Value = 5521.821599999999
$ws.Cells.Item(134, 13).Value = -2986.821599999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1800.625
$ws.Cells.Item(16, 9).Value = 2235.8572
$ws.Cells.Item(16, 10).Value = 1462.1111
$ws.Cells.Item(16, 11).Value = 2235.8572
$ws.Cells.Item(16, 12).Value = 1462.1111
$ws.Cells.Item(16, 13).Value = -1948.8572
$ws.Cells.Item(16, 14).Value = -2036.1111

$ws.Cells.Item(22, 8).Value = 440.62964
$ws.Cells.Item(22, 9).Value = 313.72726
$ws.Cells.Item(22, 11).Value = 313.72726
$ws.Cells.Item(22, 13).Value = 36.27274

$ws.Cells.Item(58, 8).Value = 1315.2307
$ws.Cells.Item(58, 9).Value = 970.8421
$ws.Cells.Item(58, 10).Value = 2250
$ws.Cells.Item(58, 11).Value = 970.8421
$ws.Cells.Item(58, 12).Value = 2250
$ws.Cells.Item(58, 13).Value = -767.8421
$ws.Cells.Item(58, 14).Value = -2656

$ws.Cells.Item(113, 8).Value = 1800.625
$ws.Cells.Item(113, 9).Value = 2235.8572
$ws.Cells.Item(113, 10).Value = 1462.1111
$ws.Cells.Item(113, 11).Value = 2235.8572
$ws.Cells.Item(113, 12).Value = 1462.1111
$ws.Cells.Item(113, 13).Value = -65.85719999999992
$ws.Cells.Item(113, 14).Value = -5802.1111

$ws.Cells.Item(132, 8).Value = 55559216
$ws.Cells.Item(132, 9).Value = 83337336
$ws.Cells.Item(132, 10).Value = 27781094
$ws.Cells.Item(132, 11).Value = 250012008
$ws.Cells.Item(132, 12).Value = 83343282
$ws.Cells.Item(132, 13).Value = -250009478
$ws.Cells.Item(132, 14).Value = -83348342

$ws.Cells.Item(134, 8).Value = 5106046
$ws.Cells.Item(134, 9).Value = 6948658.5
$ws.Cells.Item(134, 10).Value = 3425.6155
$ws.Cells.Item(134, 11).Value = 20845975.5
$ws.Cells.Item(134, 12).Value = 10276.8465
$ws.Cells.Item(134, 13).Value = -20843440.5
$ws.Cells.Item(134, 14).Value = -15346.8465

$ws.Cells.Item(136, 8).Value = 1315.2307
$ws.Cells.Item(136, 9).Value = 970.8421
$ws.Cells.Item(136, 10).Value = 2250
$ws.Cells.Item(136, 11).Value = 2912.5263
$ws.Cells.Item(136, 12).Value = 6750
$ws.Cells.Item(136, 13).Value = -362.5263
$ws.Cells.Item(136, 14).Value = -11850


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1052.4791
$ws.Cells.Item(5, 9).Value = 729
$ws.Cells.Item(5, 10).Value = 1764.1333
$ws.Cells.Item(5, 11).Value = 2187
$ws.Cells.Item(5, 12).Value = 5292.3999
$ws.Cells.Item(5, 13).Value = -2075
$ws.Cells.Item(5, 14).Value = -5516.3999

$ws.Cells.Item(58, 8).Value = 1571.2963
$ws.Cells.Item(58, 9).Value = 713.8889
$ws.Cells.Item(58, 11).Value = 2141.6667
$ws.Cells.Item(58, 13).Value = -2013.6667

$ws.Cells.Item(100, 8).Value = 8828
$ws.Cells.Item(100, 10).Value = 8828
$ws.Cells.Item(100, 12).Value = 26484
$ws.Cells.Item(100, 14).Value = -28106

$ws.Cells.Item(106, 8).Value = 9743.25
$ws.Cells.Item(106, 10).Value = 9743.25
$ws.Cells.Item(106, 12).Value = 29229.75
$ws.Cells.Item(106, 14).Value = -31121.75

$ws.Cells.Item(135, 8).Value = 1052.4791
$ws.Cells.Item(135, 9).Value = 729
$ws.Cells.Item(135, 10).Value = 1764.1333
$ws.Cells.Item(135, 11).Value = 6561
$ws.Cells.Item(135, 12).Value = 15877.1997
$ws.Cells.Item(135, 13).Value = -4026
$ws.Cells.Item(135, 14).Value = -20947.1997


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 80546.14
$ws.Cells.Item(113, 9).Value = 101812.09
$ws.Cells.Item(113, 10).Value = 2571
$ws.Cells.Item(113, 11).Value = 101812.09
$ws.Cells.Item(113, 12).Value = 2571
$ws.Cells.Item(113, 13).Value = -99642.09
$ws.Cells.Item(113, 14).Value = -6911

$ws.Cells.Item(125, 8).Value = 40326
$ws.Cells.Item(125, 10).Value = 40326
$ws.Cells.Item(125, 12).Value = 40326
$ws.Cells.Item(125, 14).Value = -45246

$ws.Cells.Item(132, 8).Value = 41673452
$ws.Cells.Item(132, 9).Value = 71437780
$ws.Cells.Item(132, 10).Value = 3400.5
$ws.Cells.Item(132, 11).Value = 214313340
$ws.Cells.Item(132, 12).Value = 10201.5
$ws.Cells.Item(132, 13).Value = -214310810
$ws.Cells.Item(132, 14).Value = -15261.5

$ws.Cells.Item(134, 8).Value = 28333.334
$ws.Cells.Item(134, 10).Value = 28333.334
$ws.Cells.Item(134, 12).Value = 85000.00199999999
$ws.Cells.Item(134, 14).Value = -90070.00199999999

$ws.Cells.Item(136, 8).Value = 38663
$ws.Cells.Item(136, 10).Value = 38663
$ws.Cells.Item(136, 12).Value = 115989
$ws.Cells.Item(136, 14).Value = -121089


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(63, 8).Value = 45085.555
$ws.Cells.Item(63, 10).Value = 45085.555
$ws.Cells.Item(63, 12).Value = 45085.555
$ws.Cells.Item(63, 14).Value = -46583.555

$ws.Cells.Item(66, 8).Value = 45085.555
$ws.Cells.Item(66, 10).Value = 45085.555
$ws.Cells.Item(66, 12).Value = 135256.665
$ws.Cells.Item(66, 14).Value = -142744.665

$ws.Cells.Item(69, 8).Value = 40000
$ws.Cells.Item(69, 10).Value = 40000
$ws.Cells.Item(69, 12).Value = 40000
$ws.Cells.Item(69, 14).Value = -41622

$ws.Cells.Item(72, 8).Value = 40000
$ws.Cells.Item(72, 10).Value = 40000
$ws.Cells.Item(72, 12).Value = 120000
$ws.Cells.Item(72, 14).Value = -128112


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 4044.6
$ws.Cells.Item(96, 9).Value = 3192.2856
$ws.Cells.Item(96, 10).Value = 6033.3335
$ws.Cells.Item(96, 11).Value = 3192.2856
$ws.Cells.Item(96, 12).Value = 6033.3335
$ws.Cells.Item(96, 13).Value = -1819.2856
$ws.Cells.Item(96, 14).Value = -8779.333500000001

$ws.Cells.Item(136, 8).Value = 4314.5
$ws.Cells.Item(136, 9).Value = 5135.643
$ws.Cells.Item(136, 10).Value = 3791.9546
$ws.Cells.Item(136, 11).Value = 15406.929
$ws.Cells.Item(136, 12).Value = 11375.8638
$ws.Cells.Item(136, 13).Value = -12856.929
$ws.Cells.Item(136, 14).Value = -16475.8638

